$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.230.29'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.419.39'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.93%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.00'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.83'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.37%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.534'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.417.87'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.98%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.11%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.20'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.27%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +7.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.846.34'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.091.96'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.417.69'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.95%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.18%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.73'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.26%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.93'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.44%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.12'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +7.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '577.16'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +16.09%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.12%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.37'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0933'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +8.17%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.73%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.55'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.10%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.82'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.16%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.66'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +8.06%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.01%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.87'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.69%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.73'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '148.45'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.90%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.69'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.62%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +12.92%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '150.75'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.65%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0543'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.35'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.55%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.95%  '
